$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3994.111
$ws.Range("I74").Value = 3994.111
$ws.Range("K74").Value = 3994.111
$ws.Range("M74").Value = -3058.111
$ws.Range("H77").Value = 3994.111
$ws.Range("I77").Value = 3994.111
$ws.Range("K77").Value = 19970.555
$ws.Range("M77").Value = -15290.555
$ws.Range("H129").Value = 1836.1666
$ws.Range("I129").Value = 879.25
$ws.Range("K129").Value = 2637.75
$ws.Range("M129").Value = 2362.25
$ws.Range("H131").Value = 6870.207
$ws.Range("I131").Value = 5486.2104
$ws.Range("K131").Value = 16458.6312
$ws.Range("M131").Value = -11418.6312
$ws.Range("H137").Value = 12068.5
$ws.Range("I137").Value = 12068.5
$ws.Range("K137").Value = 36205.5
$ws.Range("M137").Value = -33655.5
$ws.Range("H141").Value = 7099.6665
$ws.Range("I141").Value = 4782.75
$ws.Range("K141").Value = 14348.25
$ws.Range("M141").Value = -9168.25

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1512
$ws.Range("I2").Value = 1124.5555
$ws.Range("J2").Value = 4999
$ws.Range("K2").Value = 1124.5555
$ws.Range("L2").Value = 4999
$ws.Range("M2").Value = -1011.5555
$ws.Range("N2").Value = -5225
$ws.Range("H32").Value = 3680.4243
$ws.Range("I32").Value = 3680.4243
$ws.Range("K32").Value = 3680.4243
$ws.Range("M32").Value = -3393.4243
$ws.Range("H61").Value = 6397.231
$ws.Range("I61").Value = 7387
$ws.Range("K61").Value = 7387
$ws.Range("M61").Value = -7175
$ws.Range("H63").Value = 5459.75
$ws.Range("I63").Value = 2673.8572
$ws.Range("J63").Value = 9360
$ws.Range("K63").Value = 2673.8572
$ws.Range("L63").Value = 9360
$ws.Range("M63").Value = -1987.8572
$ws.Range("N63").Value = -10732
$ws.Range("H66").Value = 5459.75
$ws.Range("I66").Value = 2673.8572
$ws.Range("J66").Value = 9360
$ws.Range("K66").Value = 13369.286
$ws.Range("L66").Value = 46800
$ws.Range("M66").Value = -9937.286
$ws.Range("N66").Value = -53664
$ws.Range("H97").Value = 555.0909
$ws.Range("I97").Value = 500.65
$ws.Range("J97").Value = 1099.5
$ws.Range("K97").Value = 500.65
$ws.Range("L97").Value = 1099.5
$ws.Range("M97").Value = -4.649999999999977
$ws.Range("N97").Value = -2091.5
$ws.Range("H116").Value = 1512
$ws.Range("I116").Value = 1124.5555
$ws.Range("J116").Value = 4999
$ws.Range("K116").Value = 1124.5555
$ws.Range("L116").Value = 4999
$ws.Range("M116").Value = 1169.4445
$ws.Range("N116").Value = -9587
$ws.Range("H122").Value = 1433.1305
$ws.Range("I122").Value = 1283.6
$ws.Range("K122").Value = 3850.8
$ws.Range("M122").Value = -1400.8
$ws.Range("H136").Value = 6397.231
$ws.Range("I136").Value = 7387
$ws.Range("K136").Value = 22161
$ws.Range("M136").Value = -19611

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1512
$ws.Range("I3").Value = 1124.5555
$ws.Range("J3").Value = 4999
$ws.Range("K3").Value = 1124.5555
$ws.Range("L3").Value = 4999
$ws.Range("M3").Value = -1010.5555
$ws.Range("N3").Value = -5227
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H82").Value = 42263.5
$ws.Range("I82").Value = 42263.5
$ws.Range("K82").Value = 42263.5
$ws.Range("M82").Value = -41880.5
$ws.Range("H85").Value = 42263.5
$ws.Range("I85").Value = 42263.5
$ws.Range("K85").Value = 42263.5
$ws.Range("M85").Value = -40937.5
$ws.Range("H94").Value = 2222.8276
$ws.Range("I94").Value = 1907.8636
$ws.Range("J94").Value = 3212.7144
$ws.Range("K94").Value = 1907.8636
$ws.Range("L94").Value = 3212.7144
$ws.Range("M94").Value = -1456.8636
$ws.Range("N94").Value = -4114.7144
$ws.Range("H134").Value = 2891.1052
$ws.Range("I134").Value = 2996.1667
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 8988.500100000001
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -6453.500100000001
$ws.Range("N134").Value = -8070

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3631.2222
$ws.Range("I31").Value = 1554.174
$ws.Range("J31").Value = 15574.25
$ws.Range("K31").Value = 1554.174
$ws.Range("L31").Value = 15574.25
$ws.Range("M31").Value = -1259.174
$ws.Range("N31").Value = -16164.25
$ws.Range("H34").Value = 3631.2222
$ws.Range("I34").Value = 1554.174
$ws.Range("J34").Value = 15574.25
$ws.Range("K34").Value = 1554.174
$ws.Range("L34").Value = 15574.25
$ws.Range("M34").Value = -1352.174
$ws.Range("N34").Value = -15978.25
$ws.Range("H53").Value = 29000
$ws.Range("J53").Value = 29000
$ws.Range("L53").Value = 29000
$ws.Range("N53").Value = -30214
$ws.Range("H99").Value = 2768.75
$ws.Range("I99").Value = 2025
$ws.Range("K99").Value = 2025
$ws.Range("M99").Value = -527
$ws.Range("H122").Value = 1532.4
$ws.Range("I122").Value = 1415.5
$ws.Range("K122").Value = 4246.5
$ws.Range("M122").Value = -1796.5
$ws.Range("H126").Value = 2768.75
$ws.Range("I126").Value = 2025
$ws.Range("K126").Value = 6075
$ws.Range("M126").Value = -3605
$ws.Range("H132").Value = 2362.697
$ws.Range("I132").Value = 2368.48
$ws.Range("J132").Value = 2344.625
$ws.Range("K132").Value = 7105.440000000001
$ws.Range("L132").Value = 7033.875
$ws.Range("M132").Value = -4575.440000000001
$ws.Range("N132").Value = -12093.875

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 53210.69
$ws.Range("I120").Value = 15749
$ws.Range("K120").Value = 47247
$ws.Range("M120").Value = -42409
$ws.Range("H130").Value = 2874.75
$ws.Range("I130").Value = 2874.75
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 8624.25
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("M130").Value = -3604.25

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H126").Value = 3240.375
$ws.Range("I126").Value = 3003.3845
$ws.Range("J126").Value = 4267.3335
$ws.Range("K126").Value = 9010.1535
$ws.Range("L126").Value = 12802.0005
$ws.Range("M126").Value = -6540.1535
$ws.Range("N126").Value = -17742.0005

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4041
$ws.Range("I132").Value = 4267.6113
$ws.Range("K132").Value = 12802.8339
$ws.Range("M132").Value = -10272.8339
$ws.Range("H136").Value = 4331.375
$ws.Range("I136").Value = 4419.8
$ws.Range("J136").Value = 3005
$ws.Range("K136").Value = 13259.4
$ws.Range("L136").Value = 3005
$ws.Range("M136").Value = -10709.4
$ws.Range("N136").Value = -14115
